$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2045454545454546
$ws.Range("C2").Value = 0.5324675324675324
$ws.Range("J2").Value = 0.01298701298701299
$ws.Range("P2").Value = 0.1282467532467532
$ws.Range("S2").Value = 0.1217532467532468
$ws.Range("B3").Value = 0.002881844380403458
$ws.Range("C3").Value = 0.03170028818443804
$ws.Range("J3").Value = 0.04034582132564841
$ws.Range("P3").Value = 0.7377521613832853
$ws.Range("S3").Value = 0.1873198847262248
$ws.Range("J4").Value = 0.06862745098039216
$ws.Range("O4").Value = 0.009803921568627451
$ws.Range("P4").Value = 0.5980392156862745
$ws.Range("S4").Value = 0.3235294117647059
$ws.Range("B6").Value = 0.08583690987124463
$ws.Range("D6").Value = 0.02360515021459228
$ws.Range("F6").Value = 0.06223175965665236
$ws.Range("J6").Value = 0.2124463519313305
$ws.Range("O6").Value = 0.03004291845493562
$ws.Range("Q6").Value = 0.1330472103004292
$ws.Range("R6").Value = 0.05150214592274678
$ws.Range("S6").Value = 0.4012875536480687
$ws.Range("B7").Value = 0.1382113821138211
$ws.Range("D7").Value = 0.02168021680216802
$ws.Range("E7").Value = 0.008130081300813009
$ws.Range("F7").Value = 0.07046070460704607
$ws.Range("J7").Value = 0.1355013550135501
$ws.Range("O7").Value = 0.03523035230352303
$ws.Range("Q7").Value = 0.1246612466124661
$ws.Range("R7").Value = 0.09214092140921409
$ws.Range("S7").Value = 0.3739837398373984
$ws.Range("B8").Value = 0.09790979097909791
$ws.Range("D8").Value = 0.0198019801980198
$ws.Range("E8").Value = 0.0011001100110011
$ws.Range("F8").Value = 0.0814081408140814
$ws.Range("J8").Value = 0.1342134213421342
$ws.Range("O8").Value = 0.0187018701870187
$ws.Range("Q8").Value = 0.1408140814081408
$ws.Range("R8").Value = 0.08580858085808581
$ws.Range("S8").Value = 0.4202420242024202
$ws.Range("B9").Value = 0.116991643454039
$ws.Range("D9").Value = 0.01671309192200557
$ws.Range("F9").Value = 0.07520891364902507
$ws.Range("J9").Value = 0.1337047353760446
$ws.Range("O9").Value = 0.01949860724233983
$ws.Range("Q9").Value = 0.1615598885793872
$ws.Range("R9").Value = 0.06963788300835655
$ws.Range("S9").Value = 0.4066852367688022
$ws.Range("B10").Value = 0.1054502369668247
$ws.Range("D10").Value = 0.02567140600315956
$ws.Range("E10").Value = 0.0007898894154818325
$ws.Range("F10").Value = 0.07622432859399685
$ws.Range("J10").Value = 0.1267772511848341
$ws.Range("O10").Value = 0.0193522906793049
$ws.Range("Q10").Value = 0.1951026856240126
$ws.Range("R10").Value = 0.07819905213270142
$ws.Range("S10").Value = 0.3724328593996841
$ws.Range("G11").Value = 0.152027027027027
$ws.Range("J11").Value = 0.08952702702702703
$ws.Range("K11").Value = 0.2010135135135135
$ws.Range("L11").Value = 0.5371621621621622
$ws.Range("S11").Value = 0.02027027027027027
$ws.Range("G12").Value = 0.7353846153846154
$ws.Range("J12").Value = 0.2123076923076923
$ws.Range("K12").Value = 0.009230769230769232
$ws.Range("L12").Value = 0.01538461538461539
$ws.Range("S12").Value = 0.02769230769230769
$ws.Range("G13").Value = 0.7391304347826086
$ws.Range("J13").Value = 0.2608695652173913
$ws.Range("F15").Value = 0.01875
$ws.Range("H15").Value = 0.1520833333333333
$ws.Range("I15").Value = 0.0625
$ws.Range("J15").Value = 0.3916666666666667
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.00625
$ws.Range("O15").Value = 0.07083333333333333
$ws.Range("S15").Value = 0.23125
$ws.Range("F16").Value = 0.02659574468085106
$ws.Range("H16").Value = 0.2021276595744681
$ws.Range("I16").Value = 0.06117021276595744
$ws.Range("J16").Value = 0.375
$ws.Range("K16").Value = 0.1090425531914894
$ws.Range("M16").Value = 0.02393617021276596
$ws.Range("N16").Value = 0.002659574468085106
$ws.Range("O16").Value = 0.06382978723404255
$ws.Range("S16").Value = 0.1356382978723404
$ws.Range("F17").Value = 0.01790281329923274
$ws.Range("H17").Value = 0.1867007672634271
$ws.Range("I17").Value = 0.08951406649616368
$ws.Range("J17").Value = 0.4168797953964195
$ws.Range("K17").Value = 0.08695652173913043
$ws.Range("M17").Value = 0.01918158567774936
$ws.Range("N17").Value = 0.001278772378516624
$ws.Range("O17").Value = 0.07033248081841433
$ws.Range("S17").Value = 0.1112531969309463
$ws.Range("F18").Value = 0.01939058171745152
$ws.Range("H18").Value = 0.221606648199446
$ws.Range("I18").Value = 0.09418282548476455
$ws.Range("J18").Value = 0.3739612188365651
$ws.Range("K18").Value = 0.1135734072022161
$ws.Range("M18").Value = 0.01939058171745152
$ws.Range("N18").Value = 0.002770083102493075
$ws.Range("O18").Value = 0.09141274238227147
$ws.Range("S18").Value = 0.06371191135734072
$ws.Range("F19").Value = 0.01347068145800317
$ws.Range("H19").Value = 0.2131537242472266
$ws.Range("I19").Value = 0.08161648177496038
$ws.Range("J19").Value = 0.383122028526149
$ws.Range("K19").Value = 0.1117274167987322
$ws.Range("M19").Value = 0.0150554675118859
$ws.Range("O19").Value = 0.07052297939778129
$ws.Range("S19").Value = 0.1117274167987322
